# Add a new "Hungary" sheet, modeled on the existing "Slovakia" sheet
# (same layout/Market-row/part-number-row structure), populated with the
# Hungary-specific data, and make it the active/selected tab.

$wb = $excel.ActiveWorkbook

$slovakia = $wb.Worksheets.Item("Slovakia")

# "Move or Copy" -> "Create a copy", dropped right after Slovakia.
# This duplicates all data/styles/merged cells/page setup and makes the
# new sheet the active one (mirrors the xr:uid / sheetId / activeTab shift
# seen in the diff).
$slovakia.Copy($null, $slovakia) | Out-Null

$hungary = $wb.Worksheets.Item($wb.Worksheets.Count)
$hungary.Name = "Hungary"

# Fill in the Hungary-specific values.
$hungary.Range("B2").Value = "Hungary Market"
$hungary.Range("B4").Value = "NGC-4308/T3594/T3619"

# The Hungary row for the part number picks up the thin-border cell style
# (matching every other market sheet's B4 cell, unlike Slovakia's which
# lacked it) - copy that formatting over from a cell that already has it.
$hungary.Range("A8").Copy() | Out-Null
$hungary.Range("B4").PasteSpecial(-4122) | Out-Null

# Reset Slovakia's leftover full-sheet selection (an artifact of the
# "copy sheet" operation) and leave Hungary tab active with B4 selected.
$slovakia.Cells.Select() | Out-Null
$hungary.Select() | Out-Null
$hungary.Range("B4").Select() | Out-Null
